# Payload-Range analysis update: refresh computed mass/weight figures in Weights.xlsx
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# GLOBAL RESULTS
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GLOBAL RESULTS")
$ws.Range("C2").Value  = 890.2719999999999
$ws.Range("C6").Value  = 24689.260675789817
$ws.Range("C7").Value  = 23657.360675789823
$ws.Range("C8").Value  = 21648.024608210842
$ws.Range("C12").Value = 19098.86889584896
$ws.Range("C13").Value = 18702.86889584896
$ws.Range("C14").Value = 11970.868895848958
$ws.Range("C15").Value = 10741.325250348957
$ws.Range("C16").Value = 11155.36925034896
$ws.Range("C20").Value = 242118.93820623416
$ws.Range("C21").Value = 231999.45607123425
$ws.Range("C22").Value = 212294.60052411078
$ws.Range("C26").Value = 187295.92265747715
$ws.Range("C27").Value = 183412.48925747714
$ws.Range("C28").Value = 117394.12145747716
$ws.Range("C29").Value = 105336.41726633458
$ws.Range("C30").Value = 109396.8018589346

# ---------------------------------------------------------------------------
# FUSELAGE
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("FUSELAGE")
$ws.Range("C6").Value  = 1474.0
$ws.Range("D6").Value  = -50.86666666666667
$ws.Range("C7").Value  = 1602.0
$ws.Range("D7").Value  = -46.6
$ws.Range("C8").Value  = 3011.0
$ws.Range("D8").Value  = 0.36666666666666664
$ws.Range("C9").Value  = 1468.0
$ws.Range("D9").Value  = -51.06666666666667
$ws.Range("C12").Value = 2290.666666666666
$ws.Range("D12").Value = -23.644444444444446

# ---------------------------------------------------------------------------
# WING
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WING")
$ws.Range("C7").Value  = 2375.0
$ws.Range("D7").Value  = 18.75
$ws.Range("C8").Value  = 1784.0
$ws.Range("D8").Value  = -10.8
$ws.Range("C9").Value  = 1284.0
$ws.Range("D9").Value  = -35.8
$ws.Range("C11").Value = 2404.0
$ws.Range("D11").Value = 20.2
$ws.Range("C12").Value = 2207.0
$ws.Range("D12").Value = 10.35
$ws.Range("C13").Value = 1764.7142857142858
$ws.Range("D13").Value = -11.764285714285709

# ---------------------------------------------------------------------------
# HORIZONTAL TAIL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("HORIZONTAL TAIL")
$ws.Range("C9").Value  = 133.0
$ws.Range("D9").Value  = -77.83333333333333
$ws.Range("C10").Value = 180.0
$ws.Range("D10").Value = -69.99999999999999

# ---------------------------------------------------------------------------
# VERTICAL TAIL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("VERTICAL TAIL")
$ws.Range("C8").Value = 228.0
$ws.Range("D8").Value = -43.0
$ws.Range("C9").Value = 279.0
$ws.Range("D9").Value = -30.25

# ---------------------------------------------------------------------------
# NACELLES
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("NACELLES")
$ws.Range("C3").Value  = 488.66666666666663
$ws.Range("D3").Value  = -81.90123456790121
$ws.Range("C10").Value = 304.0
$ws.Range("D10").Value = -32.44444444444443
$ws.Range("C12").Value = 244.33333333333331
$ws.Range("C17").Value = 304.0
$ws.Range("D17").Value = -32.44444444444443
$ws.Range("C19").Value = 244.33333333333331

# ---------------------------------------------------------------------------
# LANDING GEARS
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LANDING GEARS")
$ws.Range("C5").Value = 815.0
$ws.Range("D5").Value = 63.0
$ws.Range("C6").Value = 994.0
$ws.Range("D6").Value = 98.8
$ws.Range("C7").Value = 1125.0
$ws.Range("D7").Value = 125.0
$ws.Range("C8").Value = 968.0
$ws.Range("D8").Value = 93.6
$ws.Range("C9").Value = 975.5
$ws.Range("D9").Value = 95.09999999999998
